# Update "想去人数" (F column) values for several rows across sheets,
# reflecting regenerated output data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6838
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 0
$ws1.Range("F6").Value = 149
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 58
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 1289
$ws1.Range("F11").Value = 11
$ws1.Range("F14").Value = 133
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 45
$ws1.Range("F19").Value = 0
$ws1.Range("F20").Value = 0
$ws1.Range("F21").Value = 0
$ws1.Range("F23").Value = 204
$ws1.Range("F24").Value = 151

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6838
$ws4.Range("F4").Value = 27
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F10").Value = 1289
$ws4.Range("F11").Value = 11
$ws4.Range("F13").Value = 397
$ws4.Range("F15").Value = 18
$ws4.Range("F16").Value = 0
$ws4.Range("F22").Value = 83
$ws4.Range("F23").Value = 0
